$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (89) with the next date and a repeat of the latest gold price text,
# mirroring how each day's new entry was appended to this tracking sheet.
$ws.Range("A89").Value2 = "14-12-2025"
$ws.Range("B89").Value2 = $ws.Range("B88").Value2

# Match formatting of the previous data row
$ws.Range("A88:B88").Copy()
$ws.Range("A89:B89").PasteSpecial(-4122)  # xlPasteFormats
